# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 271 }

$ws.Range("C2:C$lastRow").Value = 45181
